$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Cells.Item(2, 4) '66.561.24'
Set-TextValue $ws.Cells.Item(2, 5) '  +1.08%  '

# Row 3
Set-TextValue $ws.Cells.Item(3, 4) '3.353.37'
Set-TextValue $ws.Cells.Item(3, 5) '  +1.66%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '1.00'
Set-TextValue $ws.Cells.Item(4, 5) '  +0.16%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '190.12'
Set-TextValue $ws.Cells.Item(5, 5) '  +4.92%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '560.25'
Set-TextValue $ws.Cells.Item(6, 5) '  +0.35%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 5) '  -0.03%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '3.347.98'
Set-TextValue $ws.Cells.Item(8, 5) '  +1.65%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.584'
Set-TextValue $ws.Cells.Item(9, 5) '  -0.96%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '0.184'
Set-TextValue $ws.Cells.Item(10, 5) '  -2.26%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.587'
Set-TextValue $ws.Cells.Item(11, 5) '  -0.36%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '47.12'
Set-TextValue $ws.Cells.Item(12, 5) '  -1.27%  '

# Row 13
Set-TextValue $ws.Cells.Item(13, 5) '  +2.20%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '8.73'
Set-TextValue $ws.Cells.Item(14, 5) '  +1.84%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '3.888.15'
Set-TextValue $ws.Cells.Item(15, 5) '  +1.73%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '602.19'
Set-TextValue $ws.Cells.Item(16, 5) '  -5.24%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '66.618.42'
Set-TextValue $ws.Cells.Item(17, 5) '  +1.22%  '

# Row 18
$ws.Cells.Item(18, 2).Value = 'WrappedEther'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Cells.Item(18, 4) '3.366.74'
Set-TextValue $ws.Cells.Item(18, 5) '  +2.11%  '

# Row 19
$ws.Cells.Item(19, 2).Value = 'Chainlink'
$ws.Cells.Item(19, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Cells.Item(19, 4) '18.07'
Set-TextValue $ws.Cells.Item(19, 5) '  +0.97%  '

# Row 20
Set-TextValue $ws.Cells.Item(20, 5) '  +1.14%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '11.10'
Set-TextValue $ws.Cells.Item(21, 5) '  -3.05%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 5) '  -0.03%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '18.40'
Set-TextValue $ws.Cells.Item(23, 5) '  +3.92%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '5.04'
Set-TextValue $ws.Cells.Item(24, 5) '  +0.14%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '100.48'
Set-TextValue $ws.Cells.Item(25, 5) '  -6.14%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 5) '  +0.26%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 5) '  +1.39%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '2.77'
Set-TextValue $ws.Cells.Item(28, 5) '  +2.76%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '9.60'
Set-TextValue $ws.Cells.Item(29, 5) '  +0.84%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 5) '  +0.03%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '30.84'
Set-TextValue $ws.Cells.Item(31, 5) '  +1.21%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '6.76'
Set-TextValue $ws.Cells.Item(32, 5) '  +6.05%  '

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '3.99'
Set-TextValue $ws.Cells.Item(33, 5) '  -0.40%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '586.76'
Set-TextValue $ws.Cells.Item(34, 5) '  +5.87%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '11.05'
Set-TextValue $ws.Cells.Item(35, 5) '  -0.40%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 5) '  -0.12%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '3.725.35'
Set-TextValue $ws.Cells.Item(37, 5) '  +0.40%  '

# Row 38
Set-TextValue $ws.Cells.Item(38, 5) '  -0.10%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '56.49'
Set-TextValue $ws.Cells.Item(39, 5) '  -1.54%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '3.59'
Set-TextValue $ws.Cells.Item(40, 5) '  +5.53%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '34.05'
Set-TextValue $ws.Cells.Item(41, 5) '  +5.75%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 5) '  -0.13%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '3.26'
Set-TextValue $ws.Cells.Item(43, 5) '  -7.32%  '

# Row 44
$ws.Cells.Item(44, 2).Value = 'Kaspa'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Cells.Item(44, 4) '0.128'
Set-TextValue $ws.Cells.Item(44, 5) '  +0.55%  '

# Row 45
$ws.Cells.Item(45, 2).Value = 'Fetch.AI'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Cells.Item(45, 4) '2.69'
Set-TextValue $ws.Cells.Item(45, 5) '  -1.74%  '

# Row 46
Set-TextValue $ws.Cells.Item(46, 4) '0.344'
Set-TextValue $ws.Cells.Item(46, 5) '  +0.33%  '

# Row 47
Set-TextValue $ws.Cells.Item(47, 4) '3.39'
Set-TextValue $ws.Cells.Item(47, 5) '  +5.00%  '

# Row 48
Set-TextValue $ws.Cells.Item(48, 4) '0.0423'
Set-TextValue $ws.Cells.Item(48, 5) '  +1.75%  '

# Row 49
Set-TextValue $ws.Cells.Item(49, 5) '  +0.04%  '

# Row 50
Set-TextValue $ws.Cells.Item(50, 5) '  -1.00%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '0.999'
Set-TextValue $ws.Cells.Item(51, 5) '  +0.09%  '
